$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 0.532
$ws.Range("C5").Value = 8
$ws.Range("C6").Value = 1
$ws.Range("C8").Value = 1
$ws.Range("A9").Value = "XC9572XL-10TQG100C"

$ws.Range("E17").Value = "https://www.digikey.com.au/en/products/detail/ftdi-future-technology-devices-international-ltd/FT245RL-REEL/1836389"
$ws.Range("A17").Value = "USB serial controller"
$ws.Range("F9").Value = "Xilinx CPLD 72 I/O"



$ws.Range("A20").ClearContents()
$ws.Range("A21").ClearContents()
$ws.Range("A22").ClearContents()
$ws.Range("A23").ClearContents()

$ws.Range("A28").Value = "SMD Resistors (Aliexpress)"
$ws.Range("A29").Value = "SMD LEDs (Aliexpress)"
$ws.Range("A30").Value = "SMD Capacitors (Aliexpress)"
$ws.Range("A31").Value = "Crystals"

$ws.Range("C13").Select()
